$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 0.292742
$ws.Range("H2").Value = 0.8782260000000001
$ws.Range("I2").Value = 0.1005821958520865
$ws.Range("J2").Value = 0.1005821958520865
$ws.Range("M2").Value = 0.163273
$ws.Range("N2").Value = 0.489819
$ws.Range("Q2").Value = 0.047796864566
$ws.Range("R2").Value = 0.4301717810940001
$ws.Range("S2").Value = 0.1005821958520865
$ws.Range("T2").Value = 0.1005821958520865

# Row 3
$ws.Range("I3").Value = 0.8949043375045497
$ws.Range("J3").Value = 0.8949043375045498
$ws.Range("M3").Value = 0.163273
$ws.Range("N3").Value = 0.489819
$ws.Range("Q3").Value = 0.425260365981
$ws.Range("R3").Value = 3.827343293829
$ws.Range("S3").Value = 0.8949043375045497
$ws.Range("T3").Value = 0.8949043375045498

# Row 4
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 0.3333333333333333
$ws.Range("G4").Value = 0.01313633333333333
$ws.Range("H4").Value = 0.039409
$ws.Range("I4").Value = 0.004513466643363867
$ws.Range("J4").Value = 0.004513466643363868
$ws.Range("M4").Value = 0.163273
$ws.Range("N4").Value = 0.489819
$ws.Range("Q4").Value = 0.002144808552333333
$ws.Range("R4").Value = 0.019303276971
$ws.Range("S4").Value = 0.004513466643363867
$ws.Range("T4").Value = 0.004513466643363868
